# feat: add 2022-Q3 data
#
# 1) Insert a new worksheet named "2022-Q3" right before the existing
#    "2022-Q2" tab. It is created as a copy of "2022-Q2" so it inherits the
#    same header/column formatting (bold header row, centered bold index
#    column A, etc.), then its values are overwritten/extended to the 6
#    fund rows for the new quarter.
# 2) Insert a new summary row at the top of the "总计" (total) sheet's data
#    with the 2022-Q3 totals, shifting the existing rows down one and
#    renumbering the index column (A).

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------------
# Step 1: create the new "2022-Q3" sheet by duplicating "2022-Q2" (so it
# inherits identical header/column formatting) and placing it right before
# the old "2022-Q2" tab.
# ---------------------------------------------------------------------------
$q2Sheet.Copy($q2Sheet)
$q3Sheet = $wb.Worksheets.Item("2022-Q2 (2)")
$q3Sheet.Name = "2022-Q3"

# The template only has 4 data rows (rows 2-5 incl. header); we need 6 data
# rows (rows 2-7). Extend the styled index column / row formatting down by
# copying the last template data row's formats into the two extra rows.
$q3Sheet.Range("A4:H4").Copy()
$q3Sheet.Range("A5:H7").PasteSpecial(-4122)

# Fund-code / text columns must stay text (leading zeros, trailing zeros in
# "15.90" etc. matter), so values that look numeric are written with a
# leading apostrophe to force text storage.
$q3Data = @(
    @("010695", "华夏磐益一年定期开放混合", "15.90", "99.95", "4.46", "0.7091", 4),
    @("159611", "广发中证全指电力ETF", "16.62", "99.58", "2.56", "0.4255", 10),
    @("159625", "嘉实国证绿色电力ETF", "3.41", "98.77", "2.97", "0.1013", 10),
    @("561700", "博时中证全指电力公用事业ETF", "1.31", "98.79", "2.54", "0.0333", 10),
    @("561560", "华泰柏瑞中证全指电力公用事业ETF", "0.90", "98.29", "2.31", "0.0208", 10),
    @("562350", "银华中证全指电力公用事业ETF", "0.46", "97.99", "2.52", "0.0116", 10)
)

for ($i = 0; $i -lt $q3Data.Count; $i++) {
    $row = $i + 2
    $rec = $q3Data[$i]
    $q3Sheet.Cells.Item($row, 1).Value = $i
    $q3Sheet.Cells.Item($row, 2).Value = "'" + $rec[0]
    $q3Sheet.Cells.Item($row, 3).Value = "'" + $rec[1]
    $q3Sheet.Cells.Item($row, 4).Value = "'" + $rec[2]
    $q3Sheet.Cells.Item($row, 5).Value = "'" + $rec[3]
    $q3Sheet.Cells.Item($row, 6).Value = "'" + $rec[4]
    $q3Sheet.Cells.Item($row, 7).Value = "'" + $rec[5]
    $q3Sheet.Cells.Item($row, 8).Value = $rec[6]
}

# ---------------------------------------------------------------------------
# Step 2: insert the 2022-Q3 summary row at the top of the "总计" sheet's
# data (row 2), pushing the rest of the table down one row, and renumber the
# index column (A) for all rows. Written directly (not via Rows.Insert, to
# avoid Excel's "copy format from row above" side effect) so every row keeps
# its original per-column styling.
# ---------------------------------------------------------------------------

# Stamp row 8 (brand new row) with the same formatting as the existing last
# row (7) before filling in values, so the index column keeps its bold/
# centered style.
$totalSheet.Range("A7:D7").Copy()
$totalSheet.Range("A8:D8").PasteSpecial(-4122)

$totalRows = @(
    @("2022-Q3", 6, 1.3),
    @("2022-Q2", 3, 0.42),
    @("2022-Q1", 3, 0.4),
    @("2021-Q4", 1, 0.02),
    @("2021-Q3", 10, 6.13),
    @("2021-Q1", 2, 0.06),
    @("2020-Q4", 2, 0.07000000000000001)
)

for ($i = 0; $i -lt $totalRows.Count; $i++) {
    $row = $i + 2
    $rec = $totalRows[$i]
    $totalSheet.Cells.Item($row, 1).Value = $i
    $totalSheet.Cells.Item($row, 2).Value = $rec[0]
    $totalSheet.Cells.Item($row, 3).Value = $rec[1]
    $totalSheet.Cells.Item($row, 4).Value = $rec[2]
}

# Keep the same sheet "activated"/selected at the end as before the edit
# (the last tab, "2020-Q4").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
